$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.370.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.600.89'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.76'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.613.21'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.336'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.060.39'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.299.60'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.56'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.596.81'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '345.17'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.34%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.13'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.15'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.24%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0744'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.38%  '
$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.85'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.87'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.18'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.11%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.95'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.47'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.45%  '
$ws.Range("B39").Value = 'SuiNetwork'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.840'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.837'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.55'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.09%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '276.85'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.83%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.600'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.76'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0963'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0523'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.950.89'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.97%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.40'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.61%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.69%  '
